$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("forerunner_255s")

# Row 101: amazon entry
$ws.Range("A101").Value = 45235
$ws.Range("A101").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B101").Value = "11:19"
$ws.Range("C101").Value = 2686
$ws.Range("D101").Value = "amazon"
$ws.Range("E101").Value = "preto"

# Row 102: mercado livre entry
$ws.Range("A102").Value = 45235
$ws.Range("A102").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B102").Value = "11:20"
$ws.Range("C102").Value = 2625
$ws.Range("D102").Value = "mercado livre"
$ws.Range("E102").Value = "preto"
